$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bcde = New-Object 'object[,]' 44,4
$g = New-Object 'object[,]' 44,1

$bcde[0,0] = 0.127881588408715
$bcde[0,1] = 0.002777888934908601
$bcde[0,2] = 3.900430680208489
$bcde[0,3] = 0.496779210170732
$g[0,0] = 4.527869367722845

$bcde[1,0] = 0.3048080303191223
$bcde[1,1] = 0.3127903958511391
$bcde[1,2] = 0.8054896365839992
$bcde[1,3] = 0.496779210170732
$g[1,0] = 1.919867272924993

$bcde[2,0] = 1.459612070389937
$bcde[2,1] = 1.667794583268128
$bcde[2,2] = 0.1575252929769615
$bcde[2,3] = 0.496779210170732
$g[2,0] = 3.781711156805759

$bcde[3,0] = 3.230985683306322
$bcde[3,1] = 1.667794583268128
$bcde[3,2] = 0.8054896365839992
$bcde[3,3] = 0.496779210170732
$g[3,0] = 6.201049113329182

$bcde[4,0] = 0.6753301551942219
$bcde[4,1] = 0.3127903958511391
$bcde[4,2] = 0.8054896365839992
$bcde[4,3] = 0.496779210170732
$g[4,0] = 2.290389397800092

$bcde[5,0] = 0.04763786555579896
$bcde[5,1] = 0.04240448674262143
$bcde[5,2] = 0.8054896365839992
$bcde[5,3] = 0.496779210170732
$g[5,0] = 1.392311199053152

$bcde[6,0] = 0.6753301551942219
$bcde[6,1] = 0.04240448674262143
$bcde[6,2] = 0.1575252929769615
$bcde[6,3] = 0.496779210170732
$g[6,0] = 1.372039145084537

$bcde[7,0] = 1.459612070389937
$bcde[7,1] = 1.667794583268128
$bcde[7,2] = 3.900430680208489
$bcde[7,3] = 0.496779210170732
$g[7,0] = 7.524616544037286

$bcde[8,0] = 0.127881588408715
$bcde[8,1] = 0.002777888934908601
$bcde[8,2] = 0.8054896365839992
$bcde[8,3] = 0.496779210170732
$g[8,0] = 1.432928324098355

$bcde[9,0] = 1.459612070389937
$bcde[9,1] = 1.667794583268128
$bcde[9,2] = 3.900430680208489
$bcde[9,3] = 0.496779210170732
$g[9,0] = 7.524616544037286

$bcde[10,0] = 3.230985683306322
$bcde[10,1] = 1.667794583268128
$bcde[10,2] = 3.900430680208489
$bcde[10,3] = 0.496779210170732
$g[10,0] = 9.295990156953671

$bcde[11,0] = 3.230985683306322
$bcde[11,1] = 1.667794583268128
$bcde[11,2] = 3.900430680208489
$bcde[11,3] = 0.496779210170732
$g[11,0] = 9.295990156953671

$bcde[12,0] = 3.230985683306322
$bcde[12,1] = 1.667794583268128
$bcde[12,2] = 3.900430680208489
$bcde[12,3] = 0.496779210170732
$g[12,0] = 9.295990156953671

$bcde[13,0] = 3.230985683306322
$bcde[13,1] = 1.667794583268128
$bcde[13,2] = 0.1575252929769615
$bcde[13,3] = 0.496779210170732
$g[13,0] = 5.553084769722144

$bcde[14,0] = 3.230985683306322
$bcde[14,1] = 1.667794583268128
$bcde[14,2] = 3.900430680208489
$bcde[14,3] = 0.496779210170732
$g[14,0] = 9.295990156953671

$bcde[15,0] = 3.230985683306322
$bcde[15,1] = 1.667794583268128
$bcde[15,2] = 26.21740644021617
$bcde[15,3] = 0.496779210170732
$g[15,0] = 31.61296591696135

$bcde[16,0] = 0.3048080303191223
$bcde[16,1] = 0.3127903958511391
$bcde[16,2] = 0.8054896365839992
$bcde[16,3] = 0.496779210170732
$g[16,0] = 1.919867272924993

$bcde[17,0] = 3.230985683306322
$bcde[17,1] = 1.667794583268128
$bcde[17,2] = 26.21740644021617
$bcde[17,3] = 0.496779210170732
$g[17,0] = 31.61296591696135

$bcde[18,0] = 1.459612070389937
$bcde[18,1] = 1.667794583268128
$bcde[18,2] = 0.8054896365839992
$bcde[18,3] = 0.496779210170732
$g[18,0] = 4.429675500412797

$bcde[19,0] = 3.230985683306322
$bcde[19,1] = 1.667794583268128
$bcde[19,2] = 0.8054896365839992
$bcde[19,3] = 0.496779210170732
$g[19,0] = 6.201049113329182

$bcde[20,0] = 3.230985683306322
$bcde[20,1] = 1.667794583268128
$bcde[20,2] = 3.900430680208489
$bcde[20,3] = 0.496779210170732
$g[20,0] = 9.295990156953671

$bcde[21,0] = 3.230985683306322
$bcde[21,1] = 1.667794583268128
$bcde[21,2] = 0.1575252929769615
$bcde[21,3] = 0.496779210170732
$g[21,0] = 5.553084769722144

$bcde[22,0] = 0.6753301551942219
$bcde[22,1] = 1.667794583268128
$bcde[22,2] = 0.1575252929769615
$bcde[22,3] = 0.496779210170732
$g[22,0] = 2.997429241610044

$bcde[23,0] = 3.230985683306322
$bcde[23,1] = 1.667794583268128
$bcde[23,2] = 3.900430680208489
$bcde[23,3] = 0.496779210170732
$g[23,0] = 9.295990156953671

$bcde[24,0] = 3.230985683306322
$bcde[24,1] = 1.667794583268128
$bcde[24,2] = 0.8054896365839992
$bcde[24,3] = 0.496779210170732
$g[24,0] = 6.201049113329182

$bcde[25,0] = 3.230985683306322
$bcde[25,1] = 1.667794583268128
$bcde[25,2] = 0.8054896365839992
$bcde[25,3] = 0.496779210170732
$g[25,0] = 6.201049113329182

$bcde[26,0] = 3.230985683306322
$bcde[26,1] = 1.667794583268128
$bcde[26,2] = 0.1575252929769615
$bcde[26,3] = 0.496779210170732
$g[26,0] = 5.553084769722144

$bcde[27,0] = 3.230985683306322
$bcde[27,1] = 1.667794583268128
$bcde[27,2] = 3.900430680208489
$bcde[27,3] = 0.496779210170732
$g[27,0] = 9.295990156953671

$bcde[28,0] = 0.127881588408715
$bcde[28,1] = 0.3127903958511391
$bcde[28,2] = 0.1575252929769615
$bcde[28,3] = 0.496779210170732
$g[28,0] = 1.094976487407548

$bcde[29,0] = 0.6753301551942219
$bcde[29,1] = 1.667794583268128
$bcde[29,2] = 0.8054896365839992
$bcde[29,3] = 0.496779210170732
$g[29,0] = 3.645393585217082

$bcde[30,0] = 3.230985683306322
$bcde[30,1] = 1.667794583268128
$bcde[30,2] = 26.21740644021617
$bcde[30,3] = 0.496779210170732
$g[30,0] = 31.61296591696135

$bcde[31,0] = 0.127881588408715
$bcde[31,1] = 0.04240448674262143
$bcde[31,2] = 0.1575252929769615
$bcde[31,3] = 0.496779210170732
$g[31,0] = 0.8245905782990299

$bcde[32,0] = 0.127881588408715
$bcde[32,1] = 0.3127903958511391
$bcde[32,2] = 0.1575252929769615
$bcde[32,3] = 0.496779210170732
$g[32,0] = 1.094976487407548

$bcde[33,0] = 1.459612070389937
$bcde[33,1] = 0.3127903958511391
$bcde[33,2] = 0.8054896365839992
$bcde[33,3] = 0.496779210170732
$g[33,0] = 3.074671312995807

$bcde[34,0] = 1.459612070389937
$bcde[34,1] = 1.667794583268128
$bcde[34,2] = 3.900430680208489
$bcde[34,3] = 0.496779210170732
$g[34,0] = 7.524616544037286

$bcde[35,0] = 3.230985683306322
$bcde[35,1] = 1.667794583268128
$bcde[35,2] = 0.1575252929769615
$bcde[35,3] = 0.496779210170732
$g[35,0] = 5.553084769722144

$bcde[36,0] = 3.230985683306322
$bcde[36,1] = 1.667794583268128
$bcde[36,2] = 0.8054896365839992
$bcde[36,3] = 0.496779210170732
$g[36,0] = 6.201049113329182

$bcde[37,0] = 3.230985683306322
$bcde[37,1] = 1.667794583268128
$bcde[37,2] = 3.900430680208489
$bcde[37,3] = 0.496779210170732
$g[37,0] = 9.295990156953671

$bcde[38,0] = 0.6753301551942219
$bcde[38,1] = 1.667794583268128
$bcde[38,2] = 26.21740644021617
$bcde[38,3] = 0.496779210170732
$g[38,0] = 29.05731038884925

$bcde[39,0] = 3.230985683306322
$bcde[39,1] = 1.667794583268128
$bcde[39,2] = 0.8054896365839992
$bcde[39,3] = 0.496779210170732
$g[39,0] = 6.201049113329182

$bcde[40,0] = 3.230985683306322
$bcde[40,1] = 1.667794583268128
$bcde[40,2] = 3.900430680208489
$bcde[40,3] = 0.496779210170732
$g[40,0] = 9.295990156953671

$bcde[41,0] = 3.230985683306322
$bcde[41,1] = 1.667794583268128
$bcde[41,2] = 0.8054896365839992
$bcde[41,3] = 0.496779210170732
$g[41,0] = 6.201049113329182

$bcde[42,0] = 3.230985683306322
$bcde[42,1] = 1.667794583268128
$bcde[42,2] = 3.900430680208489
$bcde[42,3] = 0.496779210170732
$g[42,0] = 9.295990156953671

$bcde[43,0] = 1.459612070389937
$bcde[43,1] = 0.3127903958511391
$bcde[43,2] = 3.900430680208489
$bcde[43,3] = 0.496779210170732
$g[43,0] = 6.169612356620297

$ws.Range("B2:E45").Value = $bcde
$ws.Range("G2:G45").Value = $g

Write-Output "Updated s_vals for kimbrel_craig (filtered save games)"
